# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-level detail) right before the
#    "总计" (totals) sheet.
# 2. Insert a new summary row at the top of "总计" for 2022-Q1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value that must be preserved as literal TEXT (not
# re-interpreted as a number by Excel's smart entry), then strip the
# quote-prefix style that this technique leaves behind so the cell ends up
# with no explicit style index (matching plain, unstyled data cells).
# ---------------------------------------------------------------------------
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# ===========================================================================
# STEP 1 — create the "2022-Q1" sheet (fund holdings detail for the quarter)
# ===========================================================================
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Match the page margins / outline settings used by the rest of the workbook
$q1.PageSetup.LeftMargin = 0.75 * 72
$q1.PageSetup.RightMargin = 0.75 * 72
$q1.PageSetup.TopMargin = 1 * 72
$q1.PageSetup.BottomMargin = 1 * 72
$q1.PageSetup.HeaderMargin = 0.5 * 72
$q1.PageSetup.FooterMargin = 0.5 * 72
$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1

# Pull the header/index-column style (bold, centered, bordered) from an
# existing sheet so the new sheet reuses the same style index.
$styleSrcSheet = $wb.Worksheets.Item("2021-Q4")

# -- header row (B1:H1) ------------------------------------------------------
$styleSrcSheet.Range("B1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$cols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Range($cols[$i] + "1").Value = $headers[$i]
}

# -- index column (A2:A7) ----------------------------------------------------
$styleSrcSheet.Range("A2").Copy()
$q1.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# -- data rows ----------------------------------------------------------------
# code, name, scale, total position, position ratio, held value, rank
$rows = @(
    @("006567", "中泰星元价值优选灵活配置混合", "44.13", "81.95", "3.77", "1.6637", 9),
    @("013776", "中泰兴为价值精选混合A",        "20.31", "85.34", "4.85", "0.9850", 5),
    @("006624", "中泰玉衡价值优选混合",          "17.75", "81.95", "4.49", "0.7970", 7),
    @("010728", "中泰兴诚价值一年持有期混合A",    "9.94",  "88.08", "4.95", "0.4920", 7),
    @("013777", "中泰兴为价值精选混合C",          "8.71",  "85.34", "4.85", "0.4224", 5),
    @("010729", "中泰兴诚价值一年持有期混合C",    "1.72",  "88.08", "4.95", "0.0851", 7)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $r + 2
    $q1.Range("A" + $row).Value = $r

    Set-TextValue $q1.Range("B" + $row) $rows[$r][0]
    Set-TextValue $q1.Range("C" + $row) $rows[$r][1]
    Set-TextValue $q1.Range("D" + $row) $rows[$r][2]
    Set-TextValue $q1.Range("E" + $row) $rows[$r][3]
    Set-TextValue $q1.Range("F" + $row) $rows[$r][4]
    Set-TextValue $q1.Range("G" + $row) $rows[$r][5]

    $q1.Range("H" + $row).Value = $rows[$r][6]
}

# ===========================================================================
# STEP 2 — prepend a 2022-Q1 row to the "总计" summary sheet
# ===========================================================================
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 4.45
$total.Range("B2:D2").Style = "Normal"

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$total.Range("A2").Value = 0

# ===========================================================================
# Restore the originally active sheet/selection
# ===========================================================================
$wb.Worksheets.Item("2020-Q4").Activate()
